# The workbook's meeting log rows got reshuffled: for several groups of
# rows sharing the same date, the "Entity/ies met" (C) and "Subject(s)" (D)
# — and sometimes "Location" (B) — values were rotated/swapped between rows.
# Row A (date) stays put; only B/C/D text is rewritten in place.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 103-105 (12/10/2021, Brussels)
$ws.Range("C103").Value = 'Vereinigung der österreichischen Industrie - Industriellenvereinigung'
$ws.Range("D103").Value = 'The Pillar of Social rights (and an action plan), the minimum wage proposal and the initiative for social dialogue.'

$ws.Range("C104").Value = 'Industriegewerkschaft Metall'

$ws.Range("C105").Value = 'Deutscher Gewerkschaftsbund'
$ws.Range("D105").Value = 'Exchange on ongoing files in social policy. '

# Rows 118-119 (13/07/2021)
$ws.Range("B118").Value = 'Brussels'
$ws.Range("C118").Value = 'Association Internationale de la Mutualité'
$ws.Range("D118").Value = ' Meeting on the Social Economy Action Plan and the future of welfare state.'

$ws.Range("B119").Value = 'Videoconference'
$ws.Range("C119").Value = 'The European Region of the International Lesbian, Gay, Bisexual, Trans and Intersex Association, European Federation of National Organisations working with the Homeless'
$ws.Range("D119").Value = ' Meeting on LGBTIQ homelessness.  '

# Rows 124-125 (26/04/2021)
$ws.Range("C124").Value = 'Bolt'

$ws.Range("C125").Value = 'Uber'

# Rows 141-142 (18/02/2021)
$ws.Range("C141").Value = 'Workday'
$ws.Range("D141").Value = 'Meeting on the Pact for Skills and other workforce development initiatives.'

$ws.Range("C142").Value = 'FoodDrinkEurope, Centre de liaison des industries transformatrices de viande de l''UE, Federación española de Industrias de Alimentación y Bebidas, European farmers, Federazione Italiana dell''Industria Alimentare, Conseil Européen des Jeunes Agriculteurs, Suedzucker AG, Nestlé S.A., Unilever, CEMA - European Agricultural Machinery Industry Association, Institut national de recherche pour l’agriculture, l’alimentation et l’environnement, COMITE EUROPEEN des FABRICANTS de SUCRE, AgriFood Lithuania DIH, European Federation of Food, Agriculture and Tourism Trade Unions'
$ws.Range("D142").Value = 'Pact for Skills roundtable with the representatives of agri-food sector.   '

# Rows 166-167 (27/11/2020)
$ws.Range("C166").Value = 'Polska Izba Budownictwa, EUROPEAN FEDERATION OF BUILDING AND WOODWORKERS, European Association of Electrical Contractors, European Construction Industry Federation, EUREC- The association of European Renewable Energy Research Centres, Knauf Energy Solutions, Universidade do Porto, European Builders Confederation AISBL'
$ws.Range("D166").Value = 'Pact for Skills roundtable with the construction sector.'

$ws.Range("C167").Value = 'European Federation of Nurses Associations'
$ws.Range("D167").Value = 'Meeting on the challenges faced by frontline workers during the Covid-19 crisis.'

# Rows 169-170 (17/11/2020)
$ws.Range("C169").Value = 'European Transport Workers'' Federation'
$ws.Range("D169").Value = 'Meeting on precarious work, atypical contracts, bogus self-employment, temporary work agencies and social dumping.'

$ws.Range("C170").Value = 'Boerenbond vzw'
$ws.Range("D170").Value = 'Meeting on social security and seasonal workers.'

# Rows 187-188 (08/09/2020)
$ws.Range("B187").Value = 'Videoonference'
$ws.Range("C187").Value = 'Bolt, Industriegewerkschaft Metall, Universität Wien, Fagligt Fælles Forbund, Katholieke Universiteit te Leuven, Just Eat Takeaway.com N.V., clickworker GmbH, Smart'
$ws.Range("D187").Value = 'Roundtable on platform work.   '

$ws.Range("B188").Value = 'Phone call '
$ws.Range("C188").Value = 'BUSINESSEUROPE'
$ws.Range("D188").Value = ' Phone call on social policies and social dialogue.'

# Rows 193-194 (06/07/2020, Luxembourg)
$ws.Range("C193").Value = 'SES S.A.'
$ws.Range("D193").Value = 'Presentation of the SES activities. '

$ws.Range("C194").Value = 'Chambre des salariés, Lëtzebuerger Chrëschtleche Gewerkschafts-Bond'
$ws.Range("D194").Value = 'Entrevue sur investissements dans les compétences et la transition numérique pour accélérer la reprise, les stratégies numériques et industrielles, les PME. '

# Rows 201-202 (16/06/2020)
$ws.Range("B201").Value = 'Videoconference'
$ws.Range("C201").Value = 'BUSINESSEUROPE'
$ws.Range("D201").Value = 'Meeting with the Social Affairs Committee of BusinessEurope on the Commission’s response to the COVID-19 crisis in the area of employment and social right.   '

$ws.Range("B202").Value = 'Videoconference '
$ws.Range("C202").Value = 'Eurochild AISBL'
$ws.Range("D202").Value = 'Opening of the Eurochild Webinar “Child Guarantee: It’s Happening!” followed by a Question and Answer session.'

# Rows 209-211 (08/06/2020)
$ws.Range("B209").Value = 'Videoconference'
$ws.Range("C209").Value = 'BUSINESSEUROPE, SMEunited aisbl, SGI Europe, EUROPEAN TRADE UNION CONFEDERATION'
$ws.Range("D209").Value = 'High-level hearing with Executive Vice-President Dombrovskis and the Social Partners on the Action plan to implement the European Pillar of Social Rights.  '

$ws.Range("B210").Value = 'Videoconference '
$ws.Range("C210").Value = 'The Adecco Group'
$ws.Range("D210").Value = 'Videoconference on the effect of the crisis on employment and the importance of skills.  '

$ws.Range("B211").Value = 'Vidéoconférence '
$ws.Range("C211").Value = 'Le Labo de l''Economie Sociale et Solidaire, Territoires zéro chômeur de longue durée'
$ws.Range("D211").Value = 'Vidéoconférence sur la lutte contre le chômage de longue durée. '

# Rows 217-218 (15/05/2020)
$ws.Range("C217").Value = 'Volkswagen Aktiengesellschaft'
$ws.Range("D217").Value = 'Videoconference on recovery and skills.  '

$ws.Range("C218").Value = 'The Institute of International and European Affairs, Brussels'
$ws.Range("D218").Value = 'IIEA online conference "The EU''s social agenda, in particular in light of the COVID-19 pandemic and the EU''s post-crisis
recovery". '

# Rows 220-221 (13/05/2020)
$ws.Range("C220").Value = 'BUSINESSEUROPE, SGI Europe, SMEunited aisbl, EUROPEAN TRADE UNION CONFEDERATION'
$ws.Range("D220").Value = 'Videoconference on the Recovery Plan.  '

$ws.Range("C221").Value = 'Assemblée des Régions d''Europe'
$ws.Range("D221").Value = 'Videoconference on the employment and social impact of Covid-19. '

# Rows 243-244 (early Feb 2020)
$ws.Range("B243").Value = 'Bruxelles'
$ws.Range("C243").Value = 'Bundesvereinigung der Deutschen Arbeitgeberverbände e.V.'
$ws.Range("D243").Value = 'Social agenda'

$ws.Range("B244").Value = 'Brussels'
$ws.Range("C244").Value = 'Google'
$ws.Range("D244").Value = 'Skills'

# Rows 246-247 (23/01/2020)
$ws.Range("C246").Value = 'Confederazione Generale Italiana del Lavoro'
$ws.Range("D246").Value = 'Minimum wage'

$ws.Range("C247").Value = 'INCO.ORG'
$ws.Range("D247").Value = 'Economie Sociale'

# Rows 248/250 (22/01/2020, Davos)
$ws.Range("C248").Value = 'PwCIL'

$ws.Range("C250").Value = 'World Economic Forum'

# Rows 252-253 (17/01/2020, Berlin)
$ws.Range("C252").Value = 'Bundesverband der Deutschen Industrie e.V.'
$ws.Range("D252").Value = 'European Industry'

$ws.Range("C253").Value = 'Deutscher Gewerkschaftsbund'
$ws.Range("D253").Value = 'European Industry'
